$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the shared-string header used by column B ("value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# 2) Remove the old single data value in B2; the value column now starts on row 3
$ws.Range("B2").ClearContents()

# 3) New year-end date series for column A (rows 2..22)
$dates = @(38717,39082,39447,39813,40178,40543,40908,41274,41639,42004,42369,42735,43100,43465,43830,44196,44561,44926,45291,45657,46022)

# 4) Year-over-year values for column B (rows 3..21)
$values = @(1.176843378132464,1.383039815128395,-0.6300631236164866,-2.845830838597474,1.163890860292871,1.237492433423526,1.175130261101254,0.09561723522806265,-0.4807826571170737,-0.2020167505668247,-0.3527529999609147,0.2478074346218495,-0.4273761665070541,-0.6791462188813879,-2.118861353231827,1.553047647471506,-0.400533798485958,0.6918534271163068,0.2112326235108375)

# Copy the date cell's format (A2) down onto the newly-used rows in column A
# so the whole series keeps the same style (border/number-format/font/alignment).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
